# Applies the phs002504 CCDI query-cell refresh described in the commit
# "Committing changes for phs002504": the SQL text stored in the
# ParticipantsTab / StudiesTab / DiagnosisTab / SamplesTab query cells (and
# the accompanying stats query) was rewritten to the newer df_study-centric
# query set, plus a refreshed sheet selection/zoom.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ParticipantsTab query (row 2) - now selects from df_study/df_participant w/ synonym id
$qParticipants = @'
SELECT DISTINCT
    prt.participant_id AS "Participant ID",
    std.study AS "Study ID",
    COALESCE(prt.sex_at_birth, '') AS "Sex",
    COALESCE(prt.race, '') AS "Race",
    COALESCE(CAST(syn.synonym_id AS INT), '') AS "Synonym Participant ID"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN 
    df_diagnosis dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_publication pub ON std.id = pub."study.id"
LEFT JOIN 
    df_follow_up flw ON prt.id = flw."participant.id"
LEFT JOIN 
    df_sequencing_file seq ON smp.id = seq."sample.id"
LEFT JOIN 
    df_study_admin adm ON std.id = adm."study.id"
LEFT JOIN 
    df_study_personnel per ON std.id = per."study.id"
LEFT JOIN 
    df_study_funding fund ON std.id = fund."study.id"
LEFT JOIN 
    df_methylation_array_file maf ON smp.id = maf."sample.id"
LEFT JOIN 
    df_synonym syn ON prt.id = syn."participant.id"
LEFT JOIN 
    df_treatment trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
WHERE 
    std.study = 'phs002504' AND prt.sex_at_birth = 'Female'
ORDER BY 
    prt.participant_id ASC
LIMIT 100;
'@
$ws.Range("B2").Value = $qParticipants

# StatQuery (row 2, col C) - studies/participants/samples/files rollup
$qStats = @'
SELECT
    COUNT(DISTINCT std.study) AS "Studies",
    COUNT(DISTINCT prt.participant_id) AS "Participants",
    COUNT(DISTINCT smp.sample_id) AS "Samples",
    (COUNT(DISTINCT seq.id) + COUNT(DISTINCT maf.id)) AS "Files"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN 
    df_diagnosis dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_publication pub ON std.id = pub."study.id"
LEFT JOIN 
    df_follow_up flw ON prt.id = flw."participant.id"
LEFT JOIN 
    df_sequencing_file seq ON smp.id = seq."sample.id"
LEFT JOIN 
    df_study_admin adm ON std.id = adm."study.id"
LEFT JOIN 
    df_study_personnel per ON std.id = per."study.id"
LEFT JOIN 
    df_study_funding fund ON std.id = fund."study.id"
LEFT JOIN 
    df_methylation_array_file maf ON smp.id = maf."sample.id"
LEFT JOIN 
    df_synonym syn ON prt.id = syn."participant.id"
LEFT JOIN 
    df_treatment trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
WHERE 
    std.study = 'phs002504' AND prt.sex_at_birth = 'Female';
'@
$ws.Range("C2").Value = $qStats

# DiagnosisTab query (row 3) - now includes disease phase & last known survival status
$qDiagnosis = @'
SELECT DISTINCT
    prt.participant_id AS "Participant ID",
    COALESCE(smp.sample_id, '') AS "Sample ID",
    std.study AS "Study ID",
    COALESCE(dgn.diagnosis, '') AS "Diagnosis",
    COALESCE(dgn.anatomic_site, '') AS "Diagnosis Anatomic Site",
    COALESCE(dgn.diagnosis_classification_system, '') AS "Diagnosis Classification System",
    COALESCE(dgn.diagnosis_basis, '') AS "Diagnosis Basis",    
    COALESCE(dgn.disease_phase, '') AS "Disease Phase",
    COALESCE(CASE WHEN dgn.age_at_diagnosis = -999 THEN 'Not Reported' ELSE dgn.age_at_diagnosis END, 0) AS "Age at Diagnosis (days)",
    COALESCE(srv.last_known_survival_status, '') AS "Last Known Survival Status"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN 
    df_diagnosis dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_publication pub ON std.id = pub."study.id"
LEFT JOIN 
    df_follow_up flw ON prt.id = flw."participant.id"
LEFT JOIN 
    df_sequencing_file seq ON smp.id = seq."sample.id"
LEFT JOIN 
    df_study_admin adm ON std.id = adm."study.id"
LEFT JOIN 
    df_study_personnel per ON std.id = per."study.id"
LEFT JOIN 
    df_study_funding fund ON std.id = fund."study.id"
LEFT JOIN 
    df_methylation_array_file maf ON smp.id = maf."sample.id"
LEFT JOIN 
    df_synonym syn ON prt.id = syn."participant.id"
LEFT JOIN 
    df_treatment trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
WHERE 
    std.study = 'phs002504' AND prt.sex_at_birth = 'Female'
ORDER BY 
    prt.participant_id ASC;
'@
$ws.Range("B3").Value = $qDiagnosis

# StudiesTab query (row 4) - now a compact study/diagnosis summary
$qStudies = @'
SELECT
    std.study_name AS "Study Name",
    std.study AS "Study ID",
    CONCAT(dgn.diagnosis, ' (', COUNT(DISTINCT dgn.diagnosis_id), ')') AS "Diagnosis (Top 5)",
    dgn.anatomic_site AS "Diagnosis Anatomic Site (Top 5)",
    COUNT(DISTINCT prt.participant_id) AS "Number of Participants"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnosis dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_publication pub ON std.id = pub."study.id"
WHERE 
    std.study = 'phs002504' AND prt.sex_at_birth IN ('Male', 'Female', 'Unknown')

'@
$ws.Range("B4").Value = $qStudies

# SamplesTab query (row 5) - now drops the per-row diagnosis columns
$qSamples = @'
SELECT DISTINCT
    smp.sample_id AS "Sample ID",
    prt.participant_id AS "Participant ID",
    std.study AS "Study ID",
    smp.anatomic_site AS "Sample Anatomic Site",
    COALESCE(CASE WHEN smp.participant_age_at_collection = -999 THEN 'Not Reported' ELSE smp.participant_age_at_collection END, 0) AS "Age at Sample Collection (days)",
    COALESCE(smp.sample_tumor_status, '') AS "Sample Tumor Status",
    COALESCE(smp.tumor_classification, '') AS "Sample Tumor Classification"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN 
    df_diagnosis dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_publication pub ON std.id = pub."study.id"
LEFT JOIN 
    df_follow_up flw ON prt.id = flw."participant.id"
LEFT JOIN 
    df_sequencing_file seq ON smp.id = seq."sample.id"
LEFT JOIN 
    df_study_admin adm ON std.id = adm."study.id"
LEFT JOIN 
    df_study_personnel per ON std.id = per."study.id"
LEFT JOIN 
    df_study_funding fund ON std.id = fund."study.id"
LEFT JOIN 
    df_methylation_array_file maf ON smp.id = maf."sample.id"
LEFT JOIN 
    df_synonym syn ON prt.id = syn."participant.id"
LEFT JOIN 
    df_treatment trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
WHERE 
    std.study = 'phs002504' AND prt.sex_at_birth = 'Female' AND smp.sample_id IS NOT NULL
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@
$ws.Range("B5").Value = $qSamples

# Row heights settle to the content-driven autofit values recorded in the
# saved workbook (row 4's query text got noticeably shorter).
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 299.25
$ws.Rows.Item(5).RowHeight = 409.5
$ws.Rows.Item(6).RowHeight = 409.5

# Refresh the saved cursor position / zoom level on the active sheet view.
$ws.Range("C2").Select()
$excel.ActiveWindow.Zoom = 140
